$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 'BMLTC Multi-Purpose Bldg and EC'
$ws.Range("C2").Value = 14.9185376869108
$ws.Range("D2").Value = 120.786768211462
$ws.Range("E2").Value = 1641
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 3282
$ws.Range("H2").Value = 98460000

# Row 3
$ws.Range("B3").Value = 'F. Mendoza Memorial Elem Sch.'
$ws.Range("C3").Value = 14.9176780529243
$ws.Range("D3").Value = 120.767878787289
$ws.Range("E3").Value = 1671
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 3342
$ws.Range("H3").Value = 100260000
$ws.Range("K3").Value = $false

# Row 4
$ws.Range("B4").Value = 'Calumpit Sports Complex'
$ws.Range("C4").Value = 14.9185209048724
$ws.Range("D4").Value = 120.767571728115
$ws.Range("E4").Value = 947
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1894
$ws.Range("H4").Value = 56820000
$ws.Range("J4").Value = $true
$ws.Range("K4").Value = $true
$ws.Range("L4").Value = 'Built'
$ws.Range("M4").Value = ""

# Row 5
$ws.Range("B5").Value = 'Gatbuca Basketball Court'
$ws.Range("C5").Value = 14.9221390531142
$ws.Range("D5").Value = 120.766774213649
$ws.Range("E5").Value = 376
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 752
$ws.Range("H5").Value = 22560000

# Row 6
$ws.Range("B6").Value = 'San Miguel Meysulao High School'
$ws.Range("C6").Value = 14.9167991010101
$ws.Range("D6").Value = 120.742941581954
$ws.Range("E6").Value = 3464
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 6928
$ws.Range("H6").Value = 207840000

# Row 7
$ws.Range("B7").Value = 'Doña Damiana Elem School'
$ws.Range("C7").Value = 14.917701586824
$ws.Range("D7").Value = 120.743048619728
$ws.Range("E7").Value = 3135
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 6270
$ws.Range("H7").Value = 188100000
$ws.Range("K7").Value = $false

# Row 8
$ws.Range("B8").Value = 'Danga Dike'
$ws.Range("C8").Value = 14.9271290793079
$ws.Range("D8").Value = 120.75016278348
$ws.Range("E8").Value = 126
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 252
$ws.Range("H8").Value = 7560000
$ws.Range("I8").Value = $false
$ws.Range("J8").Value = $false
$ws.Range("K8").Value = $false

# Row 9
$ws.Range("B9").Value = 'Meysulao Multipurpose/E.C.'
$ws.Range("C9").Value = 14.905513184046
$ws.Range("D9").Value = 120.739161033176
$ws.Range("E9").Value = 100
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 200
$ws.Range("H9").Value = 6000000

# Row 10
$ws.Range("B10").Value = 'Calizon Dike'
$ws.Range("C10").Value = 14.9136800407707
$ws.Range("D10").Value = 120.755871075221
$ws.Range("E10").Value = 126
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 252
$ws.Range("H10").Value = 7560000
$ws.Range("I10").Value = $false
$ws.Range("J10").Value = $false
$ws.Range("K10").Value = $false

# Row 11
$ws.Range("B11").Value = 'San Marcos Elem. Sch.'
$ws.Range("C11").Value = 14.8978852342351
$ws.Range("D11").Value = 120.778807101888
$ws.Range("E11").Value = 1147
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 2294
$ws.Range("H11").Value = 68820000
$ws.Range("K11").Value = $false

# Row 12
$ws.Range("B12").Value = 'San Marcos National H.S.'
$ws.Range("C12").Value = 14.8933901983676
$ws.Range("D12").Value = 120.777840587943
$ws.Range("E12").Value = 6353
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 12706
$ws.Range("H12").Value = 381180000
$ws.Range("K12").Value = $false

# Row 13
$ws.Range("B13").Value = 'GMA Kapuso E.C.'
$ws.Range("C13").Value = 14.8930053477281
$ws.Range("D13").Value = 120.799658008805
$ws.Range("E13").Value = 200
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 400
$ws.Range("H13").Value = 12000000

# Row 14
$ws.Range("B14").Value = 'NV9 Multi-Purpose'
$ws.Range("C14").Value = 14.8980360749457
$ws.Range("D14").Value = 120.764234222054
$ws.Range("E14").Value = 2513
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 5026
$ws.Range("H14").Value = 150780000

# Row 15
$ws.Range("B15").Value = 'Frances E.C.'
$ws.Range("C15").Value = 14.9194611702998
$ws.Range("D15").Value = 120.762172685224
$ws.Range("E15").Value = 150
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 300
$ws.Range("H15").Value = 9000000

# Row 16
$ws.Range("B16").Value = 'Balungao E.C.'
$ws.Range("C16").Value = 14.9148551401837
$ws.Range("D16").Value = 120.761492937455
$ws.Range("E16").Value = 216
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 432
$ws.Range("H16").Value = 12960000

# Row 17
$ws.Range("B17").Value = 'Gugo E.C.'
$ws.Range("C17").Value = 14.9013032539823
$ws.Range("D17").Value = 120.754718122165
$ws.Range("E17").Value = 336
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 672
$ws.Range("H17").Value = 20160000

# Row 18
$ws.Range("B18").Value = 'San Marcos E.C.'
$ws.Range("C18").Value = 14.8976681596682
$ws.Range("D18").Value = 120.779685551006
$ws.Range("E18").Value = 40
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 80
$ws.Range("H18").Value = 2400000

# Row 19
$ws.Range("B19").Value = 'San Jose E.C.'
$ws.Range("C19").Value = 14.8832975084056
$ws.Range("D19").Value = 120.734533125581
$ws.Range("E19").Value = 268
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 536
$ws.Range("H19").Value = 16080000

# Row 20
$ws.Range("B20").Value = 'MEYTO Multi-Purpose/ E.C.'
$ws.Range("C20").Value = 14.8833047964844
$ws.Range("D20").Value = 120.72904705443
$ws.Range("E20").Value = 142
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 284
$ws.Range("H20").Value = 8520000

# Row 21
$ws.Range("B21").Value = 'Barangay Hall Bulusan'
$ws.Range("C21").Value = 14.9087960788938
$ws.Range("D21").Value = 120.742261855893
$ws.Range("E21").Value = 700
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 1400
$ws.Range("H21").Value = 42000000

# Row 22
$ws.Range("B22").Value = 'Brgy. Hall Sta. Lucia'
$ws.Range("C22").Value = 14.8995538574762
$ws.Range("D22").Value = 120.737428917917
$ws.Range("E22").Value = 150
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 300
$ws.Range("H22").Value = 9000000

# Row 23
$ws.Range("B23").Value = 'Mun. Covered Court'
$ws.Range("C23").Value = 14.9141384055205
$ws.Range("D23").Value = 120.764788274642
$ws.Range("E23").Value = 713
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 1426
$ws.Range("H23").Value = 42780000
